$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 131, shifting existing rows 131..211 down to 132..212
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with the new data record
$ws.Range("A131").Value = 4
$ws.Range("B131").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C131").Value = "Los Lagos"
$ws.Range("D131").Value = 44438
$ws.Range("E131").Value = 10
$ws.Range("F131").Value = "Fruta"
$ws.Range("G131").Value = 100102
$ws.Range("H131").Value = "Cítricos"
$ws.Range("I131").Value = 100102005
$ws.Range("J131").Value = "Naranja"
$ws.Range("K131").Value = "Navel Late"
$ws.Range("L131").Value = "Primera"
$ws.Range("M131").Value = 200
$ws.Range("N131").Value = 13000
$ws.Range("O131").Value = 13000
$ws.Range("P131").Value = 13000
$ws.Range("Q131").Value = "$/caja 15 kilos empedrada"
$ws.Range("R131").Value = "Región de O'Higgins"
$ws.Range("S131").Value = 867
$ws.Range("T131").Value = 15
